$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename shared string "Neutrophils" -> "Resolving-Mac" wherever it appears (Target cluster column D)
$ws.Cells.Item(6, 4).Value = "Resolving-Mac"
$ws.Cells.Item(11, 4).Value = "Resolving-Mac"

# Row 2 (FAPs -> ECs)
$ws.Cells.Item(2, 9).Value = 0.977669497583861
$ws.Cells.Item(2, 10).Value = 0.977669497583861
$ws.Cells.Item(2, 13).Value = 2.157506
$ws.Cells.Item(2, 14).Value = 6.472517999999999
$ws.Cells.Item(2, 15).Value = 0.3549648016839517
$ws.Cells.Item(2, 16).Value = 0.3549648016839516
$ws.Cells.Item(2, 17).Value = 16.21321530127
$ws.Cells.Item(2, 18).Value = 145.91893771143
$ws.Cells.Item(2, 19).Value = 0.3470382593223039
$ws.Cells.Item(2, 20).Value = 0.3470382593223038

# Row 3 (FAPs -> FAPs)
$ws.Cells.Item(3, 9).Value = 0.977669497583861
$ws.Cells.Item(3, 10).Value = 0.977669497583861
$ws.Cells.Item(3, 15).Value = 0.4793705560628122
$ws.Cells.Item(3, 16).Value = 0.4793705560628121
$ws.Cells.Item(3, 19).Value = 0.4686659707024257
$ws.Cells.Item(3, 20).Value = 0.4686659707024256

# Row 4 (FAPs -> Inflammatory-Mac)
$ws.Cells.Item(4, 9).Value = 0.977669497583861
$ws.Cells.Item(4, 10).Value = 0.977669497583861
$ws.Cells.Item(4, 13).Value = 0.018986
$ws.Cells.Item(4, 14).Value = 0.05695799999999999
$ws.Cells.Item(4, 15).Value = 0.003123681567871193
$ws.Cells.Item(4, 16).Value = 0.003123681567871192
$ws.Cells.Item(4, 17).Value = 0.14267589787
$ws.Cells.Item(4, 18).Value = 1.28408308083
$ws.Cells.Item(4, 19).Value = 0.003053928189072597
$ws.Cells.Item(4, 20).Value = 0.003053928189072596

# Row 5 (FAPs -> MuSCs)
$ws.Cells.Item(5, 9).Value = 0.977669497583861
$ws.Cells.Item(5, 10).Value = 0.977669497583861
$ws.Cells.Item(5, 13).Value = 0.9848966666666668
$ws.Cells.Item(5, 14).Value = 2.95469
$ws.Cells.Item(5, 15).Value = 0.1620406385718132
$ws.Cells.Item(5, 16).Value = 0.1620406385718132
$ws.Cells.Item(5, 17).Value = 7.401296546183334
$ws.Cells.Item(5, 18).Value = 66.61166891565
$ws.Cells.Item(5, 19).Value = 0.1584221897006726
$ws.Cells.Item(5, 20).Value = 0.1584221897006726

# Row 6 (FAPs -> Resolving-Mac, formerly Neutrophils)
$ws.Cells.Item(6, 9).Value = 0.977669497583861
$ws.Cells.Item(6, 10).Value = 0.977669497583861
$ws.Cells.Item(6, 13).Value = 0.003041
$ws.Cells.Item(6, 14).Value = 0.009122999999999999
$ws.Cells.Item(6, 15).Value = 0.0005003221135518961
$ws.Cells.Item(6, 16).Value = 0.000500322113551896
$ws.Cells.Item(6, 17).Value = 0.022852491595
$ws.Cells.Item(6, 18).Value = 0.205672424355
$ws.Cells.Item(6, 19).Value = 0.0004891496693863776
$ws.Cells.Item(6, 20).Value = 0.0004891496693863776

# Row 7 (MuSCs -> ECs)
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.171642
$ws.Cells.Item(7, 8).Value = 0.514926
$ws.Cells.Item(7, 9).Value = 0.02233050241613897
$ws.Cells.Item(7, 10).Value = 0.02233050241613898
$ws.Cells.Item(7, 13).Value = 2.157506
$ws.Cells.Item(7, 14).Value = 6.472517999999999
$ws.Cells.Item(7, 15).Value = 0.3549648016839517
$ws.Cells.Item(7, 16).Value = 0.3549648016839516
$ws.Cells.Item(7, 17).Value = 0.3703186448519999
$ws.Cells.Item(7, 18).Value = 3.332867803668
$ws.Cells.Item(7, 19).Value = 0.007926542361647774
$ws.Cells.Item(7, 20).Value = 0.007926542361647774

# Row 8 (MuSCs -> FAPs)
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.171642
$ws.Cells.Item(8, 8).Value = 0.514926
$ws.Cells.Item(8, 9).Value = 0.02233050241613897
$ws.Cells.Item(8, 10).Value = 0.02233050241613898
$ws.Cells.Item(8, 15).Value = 0.4793705560628122
$ws.Cells.Item(8, 16).Value = 0.4793705560628121
$ws.Cells.Item(8, 17).Value = 0.5001055142959999
$ws.Cells.Item(8, 18).Value = 4.500949628664
$ws.Cells.Item(8, 19).Value = 0.01070458536038651
$ws.Cells.Item(8, 20).Value = 0.01070458536038651

# Row 9 (MuSCs -> Inflammatory-Mac)
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.171642
$ws.Cells.Item(9, 8).Value = 0.514926
$ws.Cells.Item(9, 9).Value = 0.02233050241613897
$ws.Cells.Item(9, 10).Value = 0.02233050241613898
$ws.Cells.Item(9, 13).Value = 0.018986
$ws.Cells.Item(9, 14).Value = 0.05695799999999999
$ws.Cells.Item(9, 15).Value = 0.003123681567871193
$ws.Cells.Item(9, 16).Value = 0.003123681567871192
$ws.Cells.Item(9, 17).Value = 0.003258795012
$ws.Cells.Item(9, 18).Value = 0.029329155108
$ws.Cells.Item(9, 19).Value = 0.00006975337879859646
$ws.Cells.Item(9, 20).Value = 0.00006975337879859645

# Row 10 (MuSCs -> MuSCs)
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.171642
$ws.Cells.Item(10, 8).Value = 0.514926
$ws.Cells.Item(10, 9).Value = 0.02233050241613897
$ws.Cells.Item(10, 10).Value = 0.02233050241613898
$ws.Cells.Item(10, 13).Value = 0.9848966666666668
$ws.Cells.Item(10, 14).Value = 2.95469
$ws.Cells.Item(10, 15).Value = 0.1620406385718132
$ws.Cells.Item(10, 16).Value = 0.1620406385718132
$ws.Cells.Item(10, 17).Value = 0.16904963366
$ws.Cells.Item(10, 18).Value = 1.52144670294
$ws.Cells.Item(10, 19).Value = 0.003618448871140577
$ws.Cells.Item(10, 20).Value = 0.003618448871140577

# Row 11 (MuSCs -> Resolving-Mac, formerly Neutrophils)
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.171642
$ws.Cells.Item(11, 8).Value = 0.514926
$ws.Cells.Item(11, 9).Value = 0.02233050241613897
$ws.Cells.Item(11, 10).Value = 0.02233050241613898
$ws.Cells.Item(11, 13).Value = 0.003041
$ws.Cells.Item(11, 14).Value = 0.009122999999999999
$ws.Cells.Item(11, 15).Value = 0.0005003221135518961
$ws.Cells.Item(11, 16).Value = 0.000500322113551896
$ws.Cells.Item(11, 17).Value = 0.000521963322
$ws.Cells.Item(11, 18).Value = 0.004697669898
$ws.Cells.Item(11, 19).Value = 0.00001117244416551837
$ws.Cells.Item(11, 20).Value = 0.00001117244416551837
